# Updates the crypto price/volume table to the latest scraped snapshot.
# Mirrors the GitHub Actions "Updated cryptos list" commit: per-row Price (D)
# and Volume(1h) (E) refreshes, plus two rank swaps (Dai<->ImmutableX at
# rows 25/26, ApeXProtocol<->LidoDAOToken at rows 36/37) that also touch
# the Coin (B) and Link (C) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "313.76") need to be
# forced to Text first, otherwise Excel auto-converts them to a floating
# point number (losing trailing zeros / exact digits, e.g. "0.550" -> 0.55).
# Resetting the style back to 'Normal' afterwards keeps the cell format
# identical to the original (no explicit numeric format is left behind).

$ws.Range('D2').Value = '41.519.30'
$ws.Range('D3').Value = '2.489.24'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.550'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.34%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.42%  '
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('E13').Value = '  -2.68%  '
$ws.Range('D14').Value = '2.868.66'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.50'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '2.472.55'
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.795'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = '41.492.83'
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.75%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.91'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.05%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('E28').Value = '  -0.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.56%  '
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('B37').Value = 'ApeXProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.50'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -11.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.32%  '
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.12%  '
$ws.Range('D44').Value = '1.992.22'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0288'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('E46').Value = '  -6.50%  '
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').Value = '2.731.96'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '70.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.179'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.98%  '
